$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = "30.104.20"; E = "  +5.72%  " }
    3 = @{ D = "1.918.43"; E = "  +2.60%  " }
    4 = @{ D = "1.001"; E = "  -0.61%  " }
    5 = @{ D = "329.94"; E = "  +4.71%  " }
    6 = @{ D = "1.001"; E = "  -0.43%  " }
    7 = @{ D = "0.5214"; E = "  +2.50%  " }
    8 = @{ D = "0.4086"; E = "  +4.74%  " }
    9 = @{ D = "0.08507"; E = "  +2.01%  " }
    10 = @{ D = "43.16"; E = "  +3.45%  " }
    11 = @{ D = "1.126"; E = "  +2.04%  " }
    12 = @{ D = "22.50"; E = "  +10.28%  " }
    13 = @{ D = "6.428"; E = "  +3.52%  " }
    14 = @{ D = "1.936.70"; E = "  +3.76%  " }
    15 = @{ D = "7.415"; E = "  +1.92%  " }
    16 = @{ D = "1.001"; E = "  -0.55%  " }
    17 = @{ D = "95.61"; E = "  +5.03%  " }
    18 = @{ D = "0.00001112"; E = "  +1.00%  " }
    19 = @{ D = "0.06727"; E = "  +0.12%  " }
    20 = @{ D = "18.31"; E = "  +3.32%  " }
    21 = @{ D = "1.001"; E = "  -0.40%  " }
    22 = @{ D = "6.012"; E = "  +1.67%  " }
    23 = @{ D = "30.082.61"; E = "  +5.54%  " }
    24 = @{ D = "11.32"; E = "  +2.00%  " }
    25 = @{ D = "2.227"; E = "  +0.98%  " }
    26 = @{ D = "2.141.50"; E = "  +2.85%  " }
    27 = @{ D = "160.47"; E = "  +0.03%  " }
    28 = @{ D = "21.10"; E = "  +2.36%  " }
    29 = @{ D = "2.450"; E = "  +1.57%  " }
    30 = @{ D = "129.59"; E = "  +2.58%  " }
    31 = @{ D = "1.075"; E = "  +3.50%  " }
    32 = @{ D = "0.1055"; E = "  +1.45%  " }
    33 = @{ D = "6.094"; E = "  +6.17%  " }
    34 = @{ D = "3.630"; E = "  +0.52%  " }
    35 = @{ D = "0.02494"; E = "  +1.73%  " }
    36 = @{ D = "0.06619"; E = "  +0.74%  " }
    37 = @{ D = "0.2207"; E = "  +1.98%  " }
    38 = @{ D = "1.232"; E = "  +4.38%  " }
    39 = @{ D = "5.186"; E = "  +3.26%  " }
    40 = @{ D = "8.923"; E = "  +0.33%  " }
    41 = @{ D = "0.6524"; E = "  +2.49%  " }
    42 = @{ D = "1.249"; E = "  +0.87%  " }
    43 = @{ D = "11.65"; E = "  +5.14%  " }
    44 = @{ D = "0.6165"; E = "  +2.83%  " }
    45 = @{ D = "13.20"; E = "  +1.17%  " }
    46 = @{ D = "3.771"; E = "  +2.48%  " }
    47 = @{ D = "2.081"; E = "  +3.87%  " }
    48 = @{ D = "1.247"; E = "  +2.63%  " }
    49 = @{ D = "124.59"; E = "  +2.14%  " }
    50 = @{ D = "1.164"; E = "  +5.94%  " }
    51 = @{ D = "79.88"; E = "  +4.77%  " }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("D$row").Style = "Normal"
    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("E$row").Style = "Normal"
}
